# Weekly update: insert a new price record as row 455, pushing the
# existing rows 455:545 down to 456:546 (Feria Lagunitas de Puerto
# Montt - Zapallo, "Paine" / "1a (guarda)" series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 455; this shifts
# rows 455:545 down to 456:546 and grows the used range to R546.
$ws.Rows.Item(455).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A455").Value = 4
$ws.Range("B455").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C455").Value = "Los Lagos"
$ws.Range("D455").Value = 45211
$ws.Range("E455").Value = 10
$ws.Range("F455").Value = 100112045
$ws.Range("G455").Value = "Zapallo"
$ws.Range("H455").Value = "Paine"
$ws.Range("I455").Value = "1a (guarda)"
$ws.Range("J455").Value = 500
$ws.Range("K455").Value = 800
$ws.Range("L455").Value = 800
$ws.Range("M455").Value = 800
$ws.Range("N455").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O455").Value = "Región de O'Higgins"
$ws.Range("P455").Value = 800
$ws.Range("Q455").Value = 1
$ws.Range("R455").Value = "Hortaliza"
